# Fruta / hortaliza, semanal
# Insert 2 new weekly price rows at the top of the data block (rows 20-21),
# pushing the existing rows 20-52 down to rows 22-54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 20; existing rows 20:52 shift down to 22:54.
$ws.Rows("20:21").Insert()

# New row 20: weekly "Primera" quality record for Hijuelas origin.
$ws.Range("A20").Value = 6
$ws.Range("B20").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D20").Value = "01/26/2023"
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100102
$ws.Range("H20").Value = "Cítricos"
$ws.Range("I20").Value = 100102006
$ws.Range("J20").Value = "Pomelo"
$ws.Range("K20").Value = "Start Ruby"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 24
$ws.Range("N20").Value = 170000
$ws.Range("O20").Value = 180000
$ws.Range("P20").Value = 175000
$ws.Range("Q20").Value = "`$/bins (350 kilos)"
$ws.Range("R20").Value = "Hijuelas"
$ws.Range("S20").Value = 500
$ws.Range("T20").Value = 350

# New row 21: weekly "Segunda" quality record for Hijuelas origin.
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D21").Value = "01/26/2023"
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100102
$ws.Range("H21").Value = "Cítricos"
$ws.Range("I21").Value = 100102006
$ws.Range("J21").Value = "Pomelo"
$ws.Range("K21").Value = "Start Ruby"
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 20
$ws.Range("N21").Value = 140000
$ws.Range("O21").Value = 140000
$ws.Range("P21").Value = 140000
$ws.Range("Q21").Value = "`$/bins (350 kilos)"
$ws.Range("R21").Value = "Hijuelas"
$ws.Range("S21").Value = 400
$ws.Range("T21").Value = 350

$ws.Range("A1").Select()
